$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'64.246.16"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.21%  "
$ws.Range("D3").Formula = "'2.618.55"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").Formula = "'592.96"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").Formula = "'151.89"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.22%  "
$ws.Range("E7").Value = "  +0.09%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Formula = "'0.115"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.31%  "
$ws.Range("E10").Value = "  +3.81%  "
$ws.Range("D11").Formula = "'5.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +2.04%  "
$ws.Range("D12").Formula = "'0.153"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.31%  "
$ws.Range("D13").Formula = "'28.46"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +3.22%  "
$ws.Range("D14").Formula = "'3.095.80"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.30%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Formula = "'0.0000171"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +13.38%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Formula = "'64.223.95"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.38%  "
$ws.Range("D17").Formula = "'2.610.08"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.15%  "
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("E19").Value = "  +2.30%  "
$ws.Range("D20").Formula = "'350.77"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.39%  "
$ws.Range("D21").Formula = "'7.14"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.52%  "
$ws.Range("E22").Value = "  +0.11%  "
$ws.Range("D23").Formula = "'67.51"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.50%  "
$ws.Range("E24").Value = "  -0.99%  "
$ws.Range("D25").Formula = "'9.28"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("D26").Formula = "'1.64"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("D27").Formula = "'8.24"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.56%  "
$ws.Range("E28").Value = "  +1.40%  "
$ws.Range("D29").Formula = "'543.78"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("D30").Formula = "'0.998"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.12%  "
$ws.Range("D31").Formula = "'0.0" + [char]0x2083 + "0908"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +7.68%  "
$ws.Range("D32").Formula = "'2.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("E33").Value = "  +3.82%  "
$ws.Range("D34").Formula = "'5.61"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +7.55%  "
$ws.Range("D35").Formula = "'6.20"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.99%  "
$ws.Range("E36").Value = "  +2.57%  "
$ws.Range("D37").Formula = "'163.72"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.47%  "
$ws.Range("D38").Formula = "'20.14"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +3.61%  "
$ws.Range("D39").Formula = "'1.99"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +3.77%  "
$ws.Range("D40").Formula = "'1.00"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.15%  "
$ws.Range("E41").Value = "  -0.03%  "
$ws.Range("D42").Formula = "'168.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.62%  "
$ws.Range("D43").Formula = "'41.55"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.65%  "
$ws.Range("E44").Value = "  +5.06%  "
$ws.Range("D45").Formula = "'23.16"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +7.93%  "
$ws.Range("D46").Formula = "'0.0596"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.78%  "
$ws.Range("E47").Value = "  +11.77%  "
$ws.Range("D48").Formula = "'0.639"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.64%  "
$ws.Range("D49").Formula = "'0.0250"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.64%  "
$ws.Range("D50").Formula = "'0.0981"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.63%  "
$ws.Range("E51").Value = "  +0.62%  "
